$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.855.79'
$ws.Range("E2").Value = '  +0.43%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.991.57'
$ws.Range("E3").Value = '  +1.88%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '378.56'
$ws.Range("E5").Value = '  +7.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.43'
$ws.Range("E6").Value = '  +0.32%  '

$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("E8").Value = '  -0.11%  '

$ws.Range("E9").Value = '  +1.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.74'
$ws.Range("E10").Value = '  +0.97%  '

$ws.Range("E11").Value = '  -0.28%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0850'
$ws.Range("E12").Value = '  +0.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.77'
$ws.Range("E13").Value = '  +0.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.459.29'
$ws.Range("E14").Value = '  +1.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.54'
$ws.Range("E15").Value = '  +1.58%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.002.77'
$ws.Range("E16").Value = '  +2.08%  '

$ws.Range("E17").Value = '  -1.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.930.40'
$ws.Range("E18").Value = '  +0.67%  '

$ws.Range("E19").Value = '  +5.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.46'
$ws.Range("E20").Value = '  +2.32%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.25'
$ws.Range("E21").Value = '  +0.76%  '

$ws.Range("E22").Value = '  +0.99%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.92'
$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '264.35'
$ws.Range("E24").Value = '  -0.07%  '

$ws.Range("E25").Value = '  +5.07%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.47'
$ws.Range("E26").Value = '  +19.87%  '

$ws.Range("E27").Value = '  -1.82%  '

$ws.Range("E28").Value = '  -3.82%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.54'
$ws.Range("E29").Value = '  +5.30%  '

$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '26.24'
$ws.Range("E30").Value = '  -0.76%  '

$ws.Range("B31").Value = 'Dai'
$ws.Range("C31").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.01%  '

$ws.Range("E32").Value = '  -2.68%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.98'
$ws.Range("E33").Value = '  -0.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '34.82'
$ws.Range("E34").Value = '  -1.75%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '51.53'
$ws.Range("E35").Value = '  +1.33%  '

$ws.Range("E36").Value = '  -3.85%  '

$ws.Range("E37").Value = '  +3.11%  '

$ws.Range("E38").Value = '  +0.29%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.11'
$ws.Range("E39").Value = '  -3.47%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.71'
$ws.Range("E40").Value = '  -4.01%  '

$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.49'
$ws.Range("E41").Value = '  +1.67%  '

$ws.Range("E42").Value = '  -1.50%  '

$ws.Range("E43").Value = '  +1.23%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '124.37'
$ws.Range("E44").Value = '  +3.19%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.40'
$ws.Range("E45").Value = '  -1.48%  '

$ws.Range("E46").Value = '  +19.16%  '

$ws.Range("E47").Value = '  -2.90%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.045.85'
$ws.Range("E48").Value = '  -2.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.26'
$ws.Range("E49").Value = '  +1.13%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0339'
$ws.Range("E51").Value = '  +6.43%  '
